$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for every data row (2..110).
# The date needs to move forward by one day (45202 -> 45203, i.e. 2023-10-03 -> 2023-10-04).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 1
    }
}
